$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5150.25
$ws.Range("I19").Value = 5150.25
$ws.Range("K19").Value = 5150.25
$ws.Range("M19").Value = -4975.25

$ws.Range("H92").Value = 517.9
$ws.Range("I92").Value = 404.33334
$ws.Range("J92").Value = 1540
$ws.Range("K92").Value = 404.33334
$ws.Range("L92").Value = 1540
$ws.Range("M92").Value = 843.66666
$ws.Range("N92").Value = -4036

$ws.Range("H96").Value = 20834568
$ws.Range("J96").Value = 1604.8334
$ws.Range("L96").Value = 4814.5002
$ws.Range("N96").Value = -7560.5002

$ws.Range("H97").Value = 1027.619
$ws.Range("J97").Value = 1027.619
$ws.Range("L97").Value = 3082.857
$ws.Range("N97").Value = -4074.857

$ws.Range("H100").Value = 2459.6667
$ws.Range("J100").Value = 2614.2856
$ws.Range("L100").Value = 2614.2856
$ws.Range("N100").Value = -3696.2856

$ws.Range("H101").Value = 551.3333
$ws.Range("I101").Value = 551.3333
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1653.9999
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -31.99990000000003
$ws.Range("N101").ClearContents()

$ws.Range("H114").Value = 37875.5
$ws.Range("J114").Value = 37875.5
$ws.Range("L114").Value = 37875.5
$ws.Range("N114").Value = -46553.5

$ws.Range("H116").Value = 3578.1
$ws.Range("I116").Value = 1535.7142
$ws.Range("J116").Value = 4677.846
$ws.Range("K116").Value = 1535.7142
$ws.Range("L116").Value = 4677.846
$ws.Range("M116").Value = 1906.2858
$ws.Range("N116").Value = -11561.846

$ws.Range("H129").Value = 815.72
$ws.Range("J129").Value = 840.4776000000001
$ws.Range("L129").Value = 2521.4328
$ws.Range("N129").Value = -12521.4328

$ws.Range("H132").Value = 5898.3335
$ws.Range("I132").Value = 5898.3335
$ws.Range("K132").Value = 17695.0005
$ws.Range("M132").Value = -15165.0005

$ws.Range("H138").Value = 1785.7301
$ws.Range("I138").Value = 494.12122
$ws.Range("J138").Value = 3206.5
$ws.Range("K138").Value = 1482.36366
$ws.Range("L138").Value = 9619.5
$ws.Range("M138").Value = 3657.63634
$ws.Range("N138").Value = -19899.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3010.739
$ws.Range("I61").Value = 2021.75
$ws.Range("J61").Value = 5271.2856
$ws.Range("K61").Value = 2021.75
$ws.Range("L61").Value = 5271.2856
$ws.Range("M61").Value = -1809.75
$ws.Range("N61").Value = -5695.2856

$ws.Range("H97").Value = 1550.7142
$ws.Range("I97").Value = 1516.6666
$ws.Range("J97").Value = 1612
$ws.Range("K97").Value = 1516.6666
$ws.Range("L97").Value = 1612
$ws.Range("M97").Value = -1020.6666
$ws.Range("N97").Value = -2604

$ws.Range("H102").Value = 2457.8572
$ws.Range("I102").Value = 1241
$ws.Range("J102").Value = 5500
$ws.Range("K102").Value = 1241
$ws.Range("L102").Value = 5500
$ws.Range("M102").Value = 381
$ws.Range("N102").Value = -8744

$ws.Range("H136").Value = 3010.739
$ws.Range("I136").Value = 2021.75
$ws.Range("J136").Value = 5271.2856
$ws.Range("K136").Value = 6065.25
$ws.Range("L136").Value = 15813.8568
$ws.Range("M136").Value = -3515.25
$ws.Range("N136").Value = -20913.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13893.833
$ws.Range("I31").Value = 22216.555
$ws.Range("J31").Value = 5571.1113
$ws.Range("K31").Value = 22216.555
$ws.Range("L31").Value = 5571.1113
$ws.Range("M31").Value = -21921.555
$ws.Range("N31").Value = -6161.1113

$ws.Range("H34").Value = 13893.833
$ws.Range("I34").Value = 22216.555
$ws.Range("J34").Value = 5571.1113
$ws.Range("K34").Value = 22216.555
$ws.Range("L34").Value = 5571.1113
$ws.Range("M34").Value = -22014.555
$ws.Range("N34").Value = -5975.1113

$ws.Range("H86").Value = 10493.737
$ws.Range("I86").Value = 3134
$ws.Range("J86").Value = 20613.375
$ws.Range("K86").Value = 3134
$ws.Range("L86").Value = 20613.375
$ws.Range("M86").Value = -2011
$ws.Range("N86").Value = -22859.375

$ws.Range("H89").Value = 10493.737
$ws.Range("I89").Value = 3134
$ws.Range("J89").Value = 20613.375
$ws.Range("K89").Value = 15670
$ws.Range("L89").Value = 103066.875
$ws.Range("M89").Value = -10054
$ws.Range("N89").Value = -114298.875

$ws.Range("H94").Value = 3238.3125
$ws.Range("I94").Value = 2194.5715
$ws.Range("J94").Value = 4050.111
$ws.Range("K94").Value = 2194.5715
$ws.Range("L94").Value = 4050.111
$ws.Range("M94").Value = -1743.5715
$ws.Range("N94").Value = -4952.111

$ws.Range("H132").Value = 30551.055
$ws.Range("I132").Value = 34460.465
$ws.Range("K132").Value = 103381.395
$ws.Range("M132").Value = -100851.395

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 8333
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 8333
$ws.Range("K110").Value = 0
$ws.Range("N110").Value = -33179
$ws.Range("L110").Value = 24999
$ws.Range("M110").ClearContents()

$ws.Range("H131").Value = 775.45
$ws.Range("J131").Value = 775.45
$ws.Range("L131").Value = 2326.35
$ws.Range("N131").Value = -12406.35

$ws.Range("H140").Value = 1431.3334
$ws.Range("I140").Value = 834.44446
$ws.Range("K140").Value = 2503.33338
$ws.Range("M140").Value = 2676.66662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 26900
$ws.Range("J39").Value = 26900
$ws.Range("L39").Value = 26900
$ws.Range("N39").Value = -27964

$ws.Range("H97").Value = 1702.125
$ws.Range("I97").Value = 945.7619
$ws.Range("K97").Value = 945.7619
$ws.Range("M97").Value = -449.7619

$ws.Range("H107").Value = 1660.875
$ws.Range("I107").Value = 266.33334
$ws.Range("J107").Value = 2497.6
$ws.Range("K107").Value = 266.33334
$ws.Range("L107").Value = 2497.6
$ws.Range("M107").Value = 1653.66666
$ws.Range("N107").Value = -6337.6

$ws.Range("H113").Value = 2539
$ws.Range("I113").Value = 1949.3077
$ws.Range("K113").Value = 1949.3077
$ws.Range("M113").Value = 220.6922999999999

$ws.Range("H126").Value = 4030.125
$ws.Range("I126").Value = 3221.875
$ws.Range("J126").Value = 4838.375
$ws.Range("K126").Value = 9665.625
$ws.Range("L126").Value = 14515.125
$ws.Range("M126").Value = -7195.625
$ws.Range("N126").Value = -19455.125

$ws.Range("H132").Value = 68676.22
$ws.Range("I132").Value = 95465.82000000001
$ws.Range("K132").Value = 286397.46
$ws.Range("M132").Value = -283867.46

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3931.125
$ws.Range("I7").Value = 4411.1113
$ws.Range("J7").Value = 3314
$ws.Range("K7").Value = 4411.1113
$ws.Range("L7").Value = 3314
$ws.Range("M7").Value = -4299.1113
$ws.Range("N7").Value = -3538

$ws.Range("H22").Value = 2269.3635
$ws.Range("I22").Value = 1650.25
$ws.Range("J22").Value = 3920.3333
$ws.Range("K22").Value = 1650.25
$ws.Range("L22").Value = 3920.3333
$ws.Range("M22").Value = -1355.25
$ws.Range("N22").Value = -4510.3333

$ws.Range("H27").Value = 2269.3635
$ws.Range("I27").Value = 1650.25
$ws.Range("J27").Value = 3920.3333
$ws.Range("K27").Value = 1650.25
$ws.Range("L27").Value = 3920.3333
$ws.Range("M27").Value = -1543.25
$ws.Range("N27").Value = -4134.3333

$ws.Range("H98").Value = 42166.5
$ws.Range("J98").Value = 42166.5
$ws.Range("L98").Value = 42166.5
$ws.Range("N98").Value = -48156.5

$ws.Range("H126").Value = 3931.125
$ws.Range("I126").Value = 4411.1113
$ws.Range("J126").Value = 3314
$ws.Range("K126").Value = 13233.3339
$ws.Range("L126").Value = 9942
$ws.Range("M126").Value = -10763.3339
$ws.Range("N126").Value = -14882

$ws.Range("H132").Value = 806272.6
$ws.Range("I132").Value = 1097990.1
$ws.Range("J132").Value = 4049.5
$ws.Range("K132").Value = 3293970.3
$ws.Range("L132").Value = 12148.5
$ws.Range("M132").Value = -3291440.3
$ws.Range("N132").Value = -17208.5

$ws.Range("H133").Value = 29800
$ws.Range("J133").Value = 29800
$ws.Range("L133").Value = 29800
$ws.Range("N133").Value = -34860

$ws.Range("H136").Value = 29347.79
$ws.Range("I136").Value = 37900.215
$ws.Range("K136").Value = 113700.645
$ws.Range("M136").Value = -111150.645

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 753.1905
$ws.Range("I100").Value = 386.5
$ws.Range("K100").Value = 773
$ws.Range("M100").Value = -232
